# Auto-generated edit script applying the cryptos.xlsx diff
# (GitHub Actions data refresh: prices + 1h volume deltas,
#  plus two adjacent-row coin reorderings: THORChain<->Kaspa, WEMIXToken<->Dai)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '73.406.18'
$ws.Range("E2").Value = '  +1.72%  '
$ws.Range("D3").Value = '4.053.52'
$ws.Range("E3").Value = '  +0.85%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.68'
$ws.Range("E5").Value = '  +11.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.28'
$ws.Range("E6").Value = '  +1.42%  '
$ws.Range("D7").Value = '4.047.87'
$ws.Range("E7").Value = '  +0.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.690'
$ws.Range("E8").Value = '  -0.82%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.760'
$ws.Range("E10").Value = '  +1.63%  '
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.57'
$ws.Range("E12").Value = '  +12.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000325'
$ws.Range("E13").Value = '  -0.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.05'
$ws.Range("E14").Value = '  +3.69%  '
$ws.Range("D15").Value = '4.698.10'
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").Value = '4.053.87'
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.27'
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.23'
$ws.Range("E18").Value = '  +3.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.75'
$ws.Range("E19").Value = '  +0.90%  '
$ws.Range("D20").Value = '73.275.33'
$ws.Range("E20").Value = '  +1.76%  '
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '444.28'
$ws.Range("E22").Value = '  +4.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.63'
$ws.Range("E23").Value = '  +9.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '97.44'
$ws.Range("E24").Value = '  -0.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.52'
$ws.Range("E25").Value = '  +1.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.50'
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.28'
$ws.Range("E27").Value = '  +20.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.45'
$ws.Range("E28").Value = '  +1.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.92'
$ws.Range("E29").Value = '  +1.47%  '
$ws.Range("E30").Value = '  +2.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '36.95'
$ws.Range("E31").Value = '  +0.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.94'
$ws.Range("E32").Value = '  +10.72%  '
$ws.Range("E33").Value = '  +4.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.64'
$ws.Range("E34").Value = '  +1.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '689.02'
$ws.Range("E35").Value = '  +1.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '48.39'
$ws.Range("E36").Value = '  +8.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '68.35'
$ws.Range("E37").Value = '  +3.98%  '
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("D39").Value = '0.0₃0878'
$ws.Range("E39").Value = '  +6.03%  '
$ws.Range("B40").Value = 'THORChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.35'
$ws.Range("E40").Value = '  +16.96%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.148'
$ws.Range("E41").Value = '  -2.75%  '
$ws.Range("E42").Value = '  -0.78%  '
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.35'
$ws.Range("E43").Value = '  +4.73%  '
$ws.Range("B44").Value = 'Dai'
$ws.Range("C44").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0496'
$ws.Range("E45").Value = '  +1.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.151'
$ws.Range("E47").Value = '  +0.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.70'
$ws.Range("E48").Value = '  +2.98%  '
$ws.Range("E49").Value = '  -2.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.52'
$ws.Range("E50").Value = '  +7.04%  '
$ws.Range("E51").Value = '  +2.26%  '
